# The deck's closing slide ("Kde koupit EduShield?") still quoted the old
# price of the EduShield prototype shield, so it is removed entirely.
$p = $ppt.ActivePresentation

$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 299) {
        $target = $candidate
        break
    }
}

if ($target -eq $null) {
    # Fall back to the last slide if the well-known SlideID could not be
    # located (e.g. the deck was already renumbered upstream).
    $target = $p.Slides.Item($p.Slides.Count)
}

$target.Delete()
